# Apply "added 4wk low sales check" changes to the forecast workbook.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2 (W10)
$ws1.Range("D2").Value = 4
$ws1.Range("H2").Value = 1.75
$ws1.Range("L2").Value = 0.97

# Row 3 (W11)
$ws1.Range("D3").Value = 3
$ws1.Range("H3").Value = 0.79
$ws1.Range("J3").Value = "Urgent"
$ws1.Range("L3").Value = 1.13

# Row 4 (W12)
$ws1.Range("D4").Value = 3
$ws1.Range("H4").Value = 0
$ws1.Range("I4").Value = "High"
$ws1.Range("L4").Value = 0.8100000000000001

# Row 5 (W13)
$ws1.Range("L5").Value = 1.17

# Row 6 (W14)
$ws1.Range("L6").Value = 1.05

# Row 7 (W15)
$ws1.Range("D7").Value = 3
$ws1.Range("L7").Value = 0.98

# Row 8 (W16)
$ws1.Range("D8").Value = 3
$ws1.Range("L8").Value = 0.85

# Row 9 (W17)
$ws1.Range("D9").Value = 3
$ws1.Range("L9").Value = 1.04

# Row 10 (W18)
$ws1.Range("D10").Value = 3
$ws1.Range("L10").Value = 1.02

# Row 11 (W19)
$ws1.Range("L11").Value = 0.91

# Row 12 (W20)
$ws1.Range("L12").Value = 1.13

# Row 13 (W21)
$ws1.Range("L13").Value = 0.98

# Row 14 (W22)
$ws1.Range("L14").Value = 0.86

# Row 15 (W23)
$ws1.Range("L15").Value = 0.84

# Row 16 (W24)
$ws1.Range("L16").Value = 0.9

# Row 17 (W25)
$ws1.Range("L17").Value = 1.19

# --- Sheet 2: "Summary" ---
# These "numbers" are stored as text in the sheet (Total Forecast / Max Forecast
# rows), so force text format before assigning to avoid Excel auto-converting
# the value to a number.
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "48"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "28"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "15"

$ws2.Range("B12").NumberFormat = "@"
$ws2.Range("B12").Value = "4"
